$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 427.8889   # H33
$ws.Cells.Item(33, 9).Value = 479.7931   # I33
$ws.Cells.Item(33, 11).Value = 479.7931   # K33
$ws.Cells.Item(33, 13).Value = -250.7931   # M33
$ws.Cells.Item(55, 8).Value = 113.8   # H55
$ws.Cells.Item(55, 9).Value = 117.25   # I55
$ws.Cells.Item(55, 11).Value = 117.25   # K55
$ws.Cells.Item(55, 13).Value = 96.75   # M55
$ws.Cells.Item(106, 8).Value = 2777.2222   # H106
$ws.Cells.Item(106, 9).Value = 2761.875   # I106
$ws.Cells.Item(106, 10).Value = 2900   # J106
$ws.Cells.Item(106, 11).Value = 2761.875   # K106
$ws.Cells.Item(106, 12).Value = 2900   # L106
$ws.Cells.Item(106, 13).Value = -2130.875   # M106
$ws.Cells.Item(106, 14).Value = -4162   # N106
$ws.Cells.Item(132, 8).Value = 246915.52   # H132
$ws.Cells.Item(132, 9).Value = 253063.25   # I132
$ws.Cells.Item(132, 11).Value = 759189.75   # K132
$ws.Cells.Item(132, 13).Value = -756659.75   # M132
$ws.Cells.Item(137, 8).Value = 1847.9423   # H137
$ws.Cells.Item(137, 9).Value = 1443.5454   # I137
$ws.Cells.Item(137, 10).Value = 2144.5   # J137
$ws.Cells.Item(137, 11).Value = 4330.6362   # K137
$ws.Cells.Item(137, 12).Value = 6433.5   # L137
$ws.Cells.Item(137, 13).Value = -1780.6362   # M137
$ws.Cells.Item(137, 14).Value = -11533.5   # N137
$ws.Cells.Item(138, 8).Value = 1175119.4   # H138
$ws.Cells.Item(138, 9).Value = 2503.5   # I138
$ws.Cells.Item(138, 10).Value = 2170066.2   # J138
$ws.Cells.Item(138, 11).Value = 7510.5   # K138
$ws.Cells.Item(138, 12).Value = 6510198.600000001   # L138
$ws.Cells.Item(138, 13).Value = -2370.5   # M138
$ws.Cells.Item(138, 14).Value = -6520478.600000001   # N138
$ws.Cells.Item(141, 8).Value = 5120.853   # H141
$ws.Cells.Item(141, 9).Value = 2242.5454   # I141
$ws.Cells.Item(141, 10).Value = 100105   # J141
$ws.Cells.Item(141, 11).Value = 6727.6362   # K141
$ws.Cells.Item(141, 12).Value = 300315   # L141
$ws.Cells.Item(141, 13).Value = -1547.6362   # M141
$ws.Cells.Item(141, 14).Value = -310675   # N141

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9749.962   # H32
$ws.Cells.Item(32, 9).Value = 4198.536   # I32
$ws.Cells.Item(32, 10).Value = 52310.89   # J32
$ws.Cells.Item(32, 11).Value = 4198.536   # K32
$ws.Cells.Item(32, 12).Value = 52310.89   # L32
$ws.Cells.Item(32, 13).Value = -3911.536   # M32
$ws.Cells.Item(32, 14).Value = -52884.89   # N32
$ws.Cells.Item(74, 8).Value = 4951.25   # H74
$ws.Cells.Item(74, 9).Value = 714.1177   # I74
$ws.Cells.Item(74, 10).Value = 15241.429   # J74
$ws.Cells.Item(74, 11).Value = 714.1177   # K74
$ws.Cells.Item(74, 12).Value = 15241.429   # L74
$ws.Cells.Item(74, 13).Value = 159.8823   # M74
$ws.Cells.Item(74, 14).Value = -16989.429   # N74
$ws.Cells.Item(77, 8).Value = 4951.25   # H77
$ws.Cells.Item(77, 9).Value = 714.1177   # I77
$ws.Cells.Item(77, 10).Value = 15241.429   # J77
$ws.Cells.Item(77, 11).Value = 3570.5885   # K77
$ws.Cells.Item(77, 12).Value = 76207.145   # L77
$ws.Cells.Item(77, 13).Value = 797.4115000000002   # M77
$ws.Cells.Item(77, 14).Value = -84943.145   # N77
$ws.Cells.Item(135, 8).Value = 87500   # H135
$ws.Cells.Item(135, 10).Value = 87500   # J135
$ws.Cells.Item(135, 12).Value = 87500   # L135
$ws.Cells.Item(135, 14).Value = -97640   # N135

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1130.875   # H20
$ws.Cells.Item(20, 9).Value = 887.9231   # I20
$ws.Cells.Item(20, 10).Value = 1418   # J20
$ws.Cells.Item(20, 11).Value = 887.9231   # K20
$ws.Cells.Item(20, 12).Value = 1418   # L20
$ws.Cells.Item(20, 13).Value = -640.9231   # M20
$ws.Cells.Item(20, 14).Value = -1912   # N20
$ws.Cells.Item(134, 8).Value = 8432.777   # H134
$ws.Cells.Item(134, 9).Value = 9292.666999999999   # I134
$ws.Cells.Item(134, 10).Value = 4133.3335   # J134
$ws.Cells.Item(134, 11).Value = 27878.001   # K134
$ws.Cells.Item(134, 12).Value = 12400.0005   # L134
$ws.Cells.Item(134, 13).Value = -25343.001   # M134
$ws.Cells.Item(134, 14).Value = -17470.0005   # N134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1428.9272   # H31
$ws.Cells.Item(31, 9).Value = 1112.75   # I31
$ws.Cells.Item(31, 11).Value = 1112.75   # K31
$ws.Cells.Item(31, 13).Value = -817.75   # M31
$ws.Cells.Item(34, 8).Value = 1428.9272   # H34
$ws.Cells.Item(34, 9).Value = 1112.75   # I34
$ws.Cells.Item(34, 11).Value = 1112.75   # K34
$ws.Cells.Item(34, 13).Value = -910.75   # M34
$ws.Cells.Item(94, 8).Value = 1725.6666   # H94
$ws.Cells.Item(94, 9).Value = 1656   # I94
$ws.Cells.Item(94, 10).Value = 1760.5   # J94
$ws.Cells.Item(94, 11).Value = 1656   # K94
$ws.Cells.Item(94, 12).Value = 1760.5   # L94
$ws.Cells.Item(94, 13).Value = -1205   # M94
$ws.Cells.Item(94, 14).Value = -2662.5   # N94
$ws.Cells.Item(134, 8).Value = 2715.4062   # H134
$ws.Cells.Item(134, 9).Value = 2938.087   # I134
$ws.Cells.Item(134, 10).Value = 2146.3333   # J134
$ws.Cells.Item(134, 11).Value = 8814.261   # K134
$ws.Cells.Item(134, 12).Value = 6438.999899999999   # L134
$ws.Cells.Item(134, 13).Value = -6279.261   # M134
$ws.Cells.Item(134, 14).Value = -11508.9999   # N134

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 1676.1538   # H3
$ws.Cells.Item(3, 9).Value = 848.75   # I3
$ws.Cells.Item(3, 11).Value = 2546.25   # K3
$ws.Cells.Item(3, 13).Value = -2434.25   # M3
$ws.Cells.Item(39, 8).Value = 3414.4285   # H39
$ws.Cells.Item(39, 9).Value = 993   # I39
$ws.Cells.Item(39, 10).Value = 3818   # J39
$ws.Cells.Item(39, 11).Value = 2979   # K39
$ws.Cells.Item(39, 12).Value = 11454   # L39
$ws.Cells.Item(39, 13).Value = -2685   # M39
$ws.Cells.Item(39, 14).Value = -12042   # N39
$ws.Cells.Item(40, 8).Value = 70.111115   # H40
$ws.Cells.Item(40, 9).Value = 70.111115   # I40
$ws.Cells.Item(40, 10).Value = 0   # J40
$ws.Cells.Item(40, 11).Value = 280.44446   # K40
$ws.Cells.Item(40, 12).Value = 0   # L40
$ws.Cells.Item(40, 13).Value = -211.44446   # M40
$ws.Cells.Item(40, 14).ClearContents()   # N40
$ws.Cells.Item(68, 8).Value = 922.4524   # H68
$ws.Cells.Item(68, 9).Value = 713.2093   # I68
$ws.Cells.Item(68, 10).Value = 1141.9025   # J68
$ws.Cells.Item(68, 11).Value = 2139.6279   # K68
$ws.Cells.Item(68, 12).Value = 3425.7075   # L68
$ws.Cells.Item(68, 13).Value = -1328.6279   # M68
$ws.Cells.Item(68, 14).Value = -5047.7075   # N68
$ws.Cells.Item(71, 8).Value = 922.4524   # H71
$ws.Cells.Item(71, 9).Value = 713.2093   # I71
$ws.Cells.Item(71, 10).Value = 1141.9025   # J71
$ws.Cells.Item(71, 11).Value = 6418.8837   # K71
$ws.Cells.Item(71, 12).Value = 10277.1225   # L71
$ws.Cells.Item(71, 13).Value = -2362.8837   # M71
$ws.Cells.Item(71, 14).Value = -18389.1225   # N71
$ws.Cells.Item(107, 8).Value = 54766.812   # H107
$ws.Cells.Item(107, 10).Value = 201180.6   # J107
$ws.Cells.Item(107, 12).Value = 603541.8   # L107
$ws.Cells.Item(107, 14).Value = -607381.8   # N107
$ws.Cells.Item(129, 8).Value = 860.2222   # H129
$ws.Cells.Item(129, 9).Value = 338.4   # I129
$ws.Cells.Item(129, 10).Value = 1512.5   # J129
$ws.Cells.Item(129, 11).Value = 1015.2   # K129
$ws.Cells.Item(129, 12).Value = 4537.5   # L129
$ws.Cells.Item(129, 13).Value = 3984.8   # M129
$ws.Cells.Item(129, 14).Value = -14537.5   # N129
$ws.Cells.Item(136, 8).Value = 1944.4546   # H136
$ws.Cells.Item(136, 9).Value = 1548.625   # I136
$ws.Cells.Item(136, 11).Value = 4645.875   # K136
$ws.Cells.Item(136, 13).Value = 454.125   # M136

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 25333.334   # H26
$ws.Cells.Item(26, 10).Value = 25333.334   # J26
$ws.Cells.Item(26, 12).Value = 25333.334   # L26
$ws.Cells.Item(26, 14).Value = -25893.334   # N26
$ws.Cells.Item(50, 8).Value = 25333.334   # H50
$ws.Cells.Item(50, 10).Value = 25333.334   # J50
$ws.Cells.Item(50, 12).Value = 25333.334   # L50
$ws.Cells.Item(50, 14).Value = -26329.334   # N50
$ws.Cells.Item(70, 8).Value = 4608.0586   # H70
$ws.Cells.Item(70, 9).Value = 4533.615   # I70
$ws.Cells.Item(70, 10).Value = 4850   # J70
$ws.Cells.Item(70, 11).Value = 4533.615   # K70
$ws.Cells.Item(70, 12).Value = 4850   # L70
$ws.Cells.Item(70, 13).Value = -4263.615   # M70
$ws.Cells.Item(70, 14).Value = -5390   # N70
$ws.Cells.Item(73, 8).Value = 4608.0586   # H73
$ws.Cells.Item(73, 9).Value = 4533.615   # I73
$ws.Cells.Item(73, 10).Value = 4850   # J73
$ws.Cells.Item(73, 11).Value = 4533.615   # K73
$ws.Cells.Item(73, 12).Value = 4850   # L73
$ws.Cells.Item(73, 13).Value = -3597.615   # M73
$ws.Cells.Item(73, 14).Value = -6722   # N73

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2126.8   # H7
$ws.Cells.Item(7, 9).Value = 2126.8   # I7
$ws.Cells.Item(7, 10).Value = 0   # J7
$ws.Cells.Item(7, 11).Value = 2126.8   # K7
$ws.Cells.Item(7, 12).Value = 0   # L7
$ws.Cells.Item(7, 13).Value = -2014.8   # M7
$ws.Cells.Item(7, 14).ClearContents()   # N7
$ws.Cells.Item(93, 8).Value = 2294.8518   # H93
$ws.Cells.Item(93, 9).Value = 1746.8   # I93
$ws.Cells.Item(93, 10).Value = 2979.9167   # J93
$ws.Cells.Item(93, 11).Value = 1746.8   # K93
$ws.Cells.Item(93, 12).Value = 2979.9167   # L93
$ws.Cells.Item(93, 13).Value = -498.8   # M93
$ws.Cells.Item(93, 14).Value = -5475.9167   # N93
$ws.Cells.Item(122, 8).Value = 15476   # H122
$ws.Cells.Item(122, 9).Value = 22221.6   # I122
$ws.Cells.Item(122, 10).Value = 4233.3335   # J122
$ws.Cells.Item(122, 11).Value = 66664.79999999999   # K122
$ws.Cells.Item(122, 12).Value = 12700.0005   # L122
$ws.Cells.Item(122, 13).Value = -64214.79999999999   # M122
$ws.Cells.Item(122, 14).Value = -17600.0005   # N122
$ws.Cells.Item(126, 8).Value = 2126.8   # H126
$ws.Cells.Item(126, 9).Value = 2126.8   # I126
$ws.Cells.Item(126, 10).Value = 0   # J126
$ws.Cells.Item(126, 11).Value = 6380.400000000001   # K126
$ws.Cells.Item(126, 12).Value = 0   # L126
$ws.Cells.Item(126, 13).Value = -3910.400000000001   # M126
$ws.Cells.Item(126, 14).ClearContents()   # N126
$ws.Cells.Item(136, 8).Value = 2123.9614   # H136
$ws.Cells.Item(136, 9).Value = 1768.85   # I136
$ws.Cells.Item(136, 10).Value = 3307.6667   # J136
$ws.Cells.Item(136, 11).Value = 5306.549999999999   # K136
$ws.Cells.Item(136, 12).Value = 9923.000100000001   # L136
$ws.Cells.Item(136, 13).Value = -2756.549999999999   # M136
$ws.Cells.Item(136, 14).Value = -15023.0001   # N136

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1164   # H136
$ws.Cells.Item(136, 9).Value = 1110.6471   # I136
$ws.Cells.Item(136, 10).Value = 1466.3334   # J136
$ws.Cells.Item(136, 11).Value = 3331.9413   # K136
$ws.Cells.Item(136, 12).Value = 4399.0002   # L136
$ws.Cells.Item(136, 13).Value = -781.9412999999995   # M136
$ws.Cells.Item(136, 14).Value = -9499.0002   # N136
